$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink bookkeeping so it doesn't end up pinned to a
# stale cell reference once the rows below it shift down.
$ws.Hyperlinks.Delete()

# Insert 6 new rows above the old "Sector Distribution Details" block (old
# row 24) to make room for a new "Number of employees / Assets / Turnover"
# breakdown table. Everything at/after row 23 shifts down by 6 rows.
$ws.Rows("23:28").Insert()

# New table header (row 23)
$ws.Range("B23").Value = "Number of employees"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"

# New table body (rows 24-27)
$ws.Range("A24").Value = "Micro"
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""

$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "1-99"
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""

$ws.Range("A26").Value = "Medium"
$ws.Range("B26").Value = "100-499"
$ws.Range("C26").Value = ""
$ws.Range("D26").Value = ""

$ws.Range("A27").Value = "Large"
$ws.Range("B27").Value = ">=500"
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""

# Re-create the hyperlink at its new location (old A48 -> new A54).
$ws.Hyperlinks.Add($ws.Range("A54"), "http://www.ic.gc.ca/eic/site/061.nsf/eng/h_02800.html")
